$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New execution timestamps (RAD re-run) for rows 2..33, column B ("Date")
$newDates = @(
    "Tue Sep 26 21:20:59 EDT 2023",
    "Tue Sep 26 21:21:12 EDT 2023",
    "Tue Sep 26 21:21:25 EDT 2023",
    "Tue Sep 26 21:21:38 EDT 2023",
    "Tue Sep 26 21:21:51 EDT 2023",
    "Tue Sep 26 21:22:04 EDT 2023",
    "Tue Sep 26 21:22:17 EDT 2023",
    "Tue Sep 26 21:22:30 EDT 2023",
    "Tue Sep 26 21:22:43 EDT 2023",
    "Tue Sep 26 21:22:56 EDT 2023",
    "Tue Sep 26 21:23:09 EDT 2023",
    "Tue Sep 26 21:23:22 EDT 2023",
    "Tue Sep 26 21:23:35 EDT 2023",
    "Tue Sep 26 21:23:48 EDT 2023",
    "Tue Sep 26 21:24:01 EDT 2023",
    "Tue Sep 26 21:24:14 EDT 2023",
    "Tue Sep 26 21:24:27 EDT 2023",
    "Tue Sep 26 21:24:40 EDT 2023",
    "Tue Sep 26 21:24:53 EDT 2023",
    "Tue Sep 26 21:25:06 EDT 2023",
    "Tue Sep 26 21:25:20 EDT 2023",
    "Tue Sep 26 21:25:33 EDT 2023",
    "Tue Sep 26 21:25:46 EDT 2023",
    "Tue Sep 26 21:25:59 EDT 2023",
    "Tue Sep 26 21:26:13 EDT 2023",
    "Tue Sep 26 21:26:26 EDT 2023",
    "Tue Sep 26 21:26:39 EDT 2023",
    "Tue Sep 26 21:26:52 EDT 2023",
    "Tue Sep 26 21:27:05 EDT 2023",
    "Tue Sep 26 21:27:18 EDT 2023",
    "Tue Sep 26 21:27:31 EDT 2023",
    "Tue Sep 26 21:27:44 EDT 2023"
)

for ($i = 0; $i -lt $newDates.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 2).Value = $newDates[$i]
}

# Row 32's Result flips from Fail to Pass
$ws.Cells.Item(32, 1).Value = "Pass"
